$d = $word.ActiveDocument

# Locate the paragraph whose text contains $searchText by scanning the
# Paragraphs collection directly (more reliable here than reusing a Range
# returned from Find for subsequent .Paragraphs/.Index lookups).
function Get-ParaIndex($searchText) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -like "*$searchText*") {
            return $i
        }
    }
    return -1
}

# --- 1. Update the letter date: "September 19, 2025" -> "September 21, 2025" ---
$dateIdx = Get-ParaIndex("September 19, 2025")
$dateRange = $d.Paragraphs.Item($dateIdx).Range
$dateRange.Text = "September 21, 2025`r"     # re-set incl. the paragraph mark so
                                              # the run keeps its formatting
$d.Paragraphs.Item($dateIdx + 1).Range.Delete()   # drop the spurious extra paragraph

# --- 2. Split the mailing address into a street line and a city/state/zip line ---
$addrIdx = Get-ParaIndex("3370 Eichers Pl, Santa Clara CA 95051")
$streetRange = $d.Paragraphs.Item($addrIdx).Range
$streetRange.Text = "3370 Eichers Pl`r"

$cityRange = $d.Paragraphs.Item($addrIdx + 1).Range
$cityRange.MoveEnd(1, -1) | Out-Null         # wdCharacter: exclude paragraph mark
$cityRange.Text = "Santa Clara, CA 95051"

# --- 3. Remove the now-superfluous blank "No Spacing" paragraph that followed
#        "... Board of Directors" ---
$bodIdx = Get-ParaIndex("Board of Directors")
$d.Paragraphs.Item($bodIdx + 1).Range.Delete()
